$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data rows 2-15: AC_NO / Addition / Deleted / Modification / Status ---
# (row 9 "141 / In Progress" split out; everything below shifts up by one,
#  row 13 Addition corrected, and a new Completed row 141 entry inserted at row 8)

$ws.Range("A2").Value = 143
$ws.Range("B2").Value = 25763
$ws.Range("C2").Value = 2123
$ws.Range("D2").Value = 4947
$ws.Range("E2").Value = "Completed"

$ws.Range("A3").Value = 139
$ws.Range("B3").Value = 25128
$ws.Range("C3").Value = 1192
$ws.Range("D3").Value = 6526
$ws.Range("E3").Value = "Completed"

$ws.Range("A4").Value = 142
$ws.Range("B4").Value = 25445
$ws.Range("C4").Value = 1557
$ws.Range("D4").Value = 5001
$ws.Range("E4").Value = "Completed"

$ws.Range("A5").Value = 144
$ws.Range("B5").Value = 42179
$ws.Range("C5").Value = 3035
$ws.Range("D5").Value = 6906
$ws.Range("E5").Value = "Completed"

$ws.Range("A6").Value = 138
$ws.Range("B6").Value = 28994
$ws.Range("C6").Value = 3714
$ws.Range("D6").Value = 6998
$ws.Range("E6").Value = "Completed"

$ws.Range("A7").Value = 140
$ws.Range("B7").Value = 9328
$ws.Range("C7").Value = 1101
$ws.Range("D7").Value = 4359
$ws.Range("E7").Value = "Completed"

$ws.Range("A8").Value = 141
$ws.Range("B8").Value = 16435
$ws.Range("C8").Value = 1581
$ws.Range("D8").Value = 3874
$ws.Range("E8").Value = "Completed"

$ws.Range("A9").Value = 146
$ws.Range("B9").Value = 36579
$ws.Range("C9").Value = 4512
$ws.Range("D9").Value = 6597
$ws.Range("E9").Value = "Completed"

$ws.Range("A10").Value = 147
$ws.Range("B10").Value = 9056
$ws.Range("C10").Value = 2982
$ws.Range("D10").Value = 4750
$ws.Range("E10").Value = "Completed"

$ws.Range("A11").Value = 148
$ws.Range("B11").Value = 19751
$ws.Range("C11").Value = 2258
$ws.Range("D11").Value = 5537
$ws.Range("E11").Value = "Completed"

$ws.Range("A12").Value = 149
$ws.Range("B12").Value = 21114
$ws.Range("C12").Value = 8203
$ws.Range("D12").Value = 9700
$ws.Range("E12").Value = "Completed"

$ws.Range("A13").Value = 150
$ws.Range("B13").Value = 16891
$ws.Range("C13").Value = 4715
$ws.Range("D13").Value = 6995
$ws.Range("E13").Value = "Completed"

$ws.Range("A14").Value = 151
$ws.Range("B14").Value = 21855
$ws.Range("C14").Value = 2516
$ws.Range("D14").Value = 5504
$ws.Range("E14").Value = "Completed"

$ws.Range("A15").Value = 188
$ws.Range("B15").Value = 34365
$ws.Range("C15").Value = 4121
$ws.Range("D15").Value = 9733
$ws.Range("E15").Value = "Completed"

# --- Row 16: drop the old record, replace with a total formula in B ---
$ws.Range("A16").ClearContents()
$ws.Range("B16").Formula = "=SUM(B2:B15)"
$ws.Range("C16").ClearContents()
$ws.Range("D16").ClearContents()
$ws.Range("E16").ClearContents()

# --- New helper total in H10 ---
$ws.Range("H10").Formula = "=1293+1006"

# --- Highlight the "Deleted" column (C2:C15) in yellow, except the newly
#     inserted row 8 (AC_NO 141) which stays unstyled ---
$ws.Range("C2:C7").Interior.Color = 65535
$ws.Range("C9:C15").Interior.Color = 65535

# --- View bits: selection ---
$ws.Range("B20").Select() | Out-Null
